{"js": "// Map of old arithmetic-problem text -> new arithmetic-problem text.\n// Each table cell in the worksheet holds exactly one \"NN\u00f7N=\" style string;\n// the edit swaps those strings in place (cell/table structure is unchanged).\nconst replacements = [\n  [\"47\u00f79=\", \"19\u00f78=\"],\n  [\"14\u00f75=\", \"25\u00f73=\"],\n  [\"23\u00f74=\", \"96\u00f79=\"],\n  [\"13\u00f72=\", \"83\u00f76=\"],\n  [\"64\u00f74=\", \"68\u00f76=\"],\n  [\"43\u00f75=\", \"12\u00f79=\"],\n  [\"76\u00f76=\", \"30\u00f77=\"],\n  [\"44\u00f75=\", \"43\u00f75=\"],\n  [\"82\u00f79=\", \"95\u00f79=\"],\n  [\"77\u00f74=\", \"59\u00f75=\"],\n  [\"12\u00f75=\", \"31\u00f74=\"],\n  [\"32\u00f77=\", \"19\u00f72=\"],\n  [\"66\u00f72=\", \"43\u00f72=\"],\n  [\"88\u00f79=\", \"54\u00f78=\"],\n  [\"99\u00f77=\", \"69\u00f76=\"],\n  [\"30\u00f78=\", \"11\u00f78=\"],\n  [\"27\u00f74=\", \"63\u00f79=\"],\n  [\"37\u00f76=\", \"32\u00f77=\"],\n  [\"60\u00f75=\", \"12\u00f76=\"],\n  [\"66\u00f78=\", \"96\u00f73=\"],\n  [\"18\u00f76=\", \"27\u00f79=\"],\n  [\"58\u00f79=\", \"41\u00f78=\"],\n  [\"52\u00f72=\", \"98\u00f78=\"],\n  [\"58\u00f72=\", \"14\u00f79=\"],\n  [\"47\u00f77=\", \"30\u00f75=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (let r = 0; r < rows.items.length; r++) {\n    const cells = rows.items[r].cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (let c = 0; c < cells.items.length; c++) {\n      const cell = table.getCell(r, c);\n      cell.load(\"value\");\n      await context.sync();\n\n      const current = cell.value;\n      const hit = replacements.find(([oldText]) => oldText === current);\n      if (hit) {\n        cell.value = hit[1];\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Map of old arithmetic-problem text -> new arithmetic-problem text.\n# Each table cell in the worksheet holds exactly one \"NN\u00f7N=\" style string;\n# the edit swaps those strings in place (cell/table structure is unchanged).\n$replacements = @{\n    \"47\u00f79=\" = \"19\u00f78=\";\n    \"14\u00f75=\" = \"25\u00f73=\";\n    \"23\u00f74=\" = \"96\u00f79=\";\n    \"13\u00f72=\" = \"83\u00f76=\";\n    \"64\u00f74=\" = \"68\u00f76=\";\n    \"43\u00f75=\" = \"12\u00f79=\";\n    \"76\u00f76=\" = \"30\u00f77=\";\n    \"44\u00f75=\" = \"43\u00f75=\";\n    \"82\u00f79=\" = \"95\u00f79=\";\n    \"77\u00f74=\" = \"59\u00f75=\";\n    \"12\u00f75=\" = \"31\u00f74=\";\n    \"32\u00f77=\" = \"19\u00f72=\";\n    \"66\u00f72=\" = \"43\u00f72=\";\n    \"88\u00f79=\" = \"54\u00f78=\";\n    \"99\u00f77=\" = \"69\u00f76=\";\n    \"30\u00f78=\" = \"11\u00f78=\";\n    \"27\u00f74=\" = \"63\u00f79=\";\n    \"37\u00f76=\" = \"32\u00f77=\";\n    \"60\u00f75=\" = \"12\u00f76=\";\n    \"66\u00f78=\" = \"96\u00f73=\";\n    \"18\u00f76=\" = \"27\u00f79=\";\n    \"58\u00f79=\" = \"41\u00f78=\";\n    \"52\u00f72=\" = \"98\u00f78=\";\n    \"58\u00f72=\" = \"14\u00f79=\";\n    \"47\u00f77=\" = \"30\u00f75=\";\n}\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    $rowCount = $t.Rows.Count\n    $colCount = $t.Columns.Count\n    for ($r = 1; $r -le $rowCount; $r++) {\n        for ($c = 1; $c -le $colCount; $c++) {\n            $cell = $null\n            try {\n                $cell = $t.Cell($r, $c)\n            } catch {\n                $cell = $null\n            }\n            if ($cell -ne $null) {\n                $current = $cell.Range.Text\n                # Cell.Range.Text carries a trailing cell-mark (CR + BEL);\n                # strip it so the lookup matches the plain problem text.\n                $current = $current.TrimEnd([char]13, [char]7)\n                if ($replacements.ContainsKey($current)) {\n                    $cell.Range.Text = $replacements[$current]\n                }\n            }\n        }\n    }\n}\n"}
